# enemyDatabase.xlsx update:
# - Insert 3 new enemy rows (Doppelganger w/ FormChange, Illyia?, Dahlia?)
#   above the existing "Doppelganger (BLUE/YELLOW/RED/???/GRAY)" rows.
# - All rows below shift down by 3 (handled automatically by the row insert).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push existing rows 21-25 (and the trailing END row) down by 3 so we can
# place the 3 new rows at 21, 22, 23.
$ws.Range("A21:A23").EntireRow.Insert()

# --- New row 21: "Doppelganger" (boss form-change variant) ---
$ws.Range("A21").Value = "Doppelganger"
$ws.Range("B21").Value = "spr_bt_doppelganger_g"
$ws.Range("C21").Value = "doppelganger_1"
$ws.Range("D21").Value = 300
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 3
$ws.Range("G21").Value = 0.75
$ws.Range("H21").Value = 0.1
$ws.Range("I21").Value = 0.5
$ws.Range("J21").Value = 1
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = "FormChange1"
$ws.Range("M21").Value = "Attacker1,Bunbuku,true"
$ws.Range("N21").Value = 0
$ws.Range("O21").Value = 0.5
$ws.Range("P21").Value = 0.5
$ws.Range("Q21").Value = 0.5
$ws.Range("R21").Value = "GROUP/DEFAULT"
$ws.Range("S21").Value = "sword"
$ws.Range("T21").Value = "null"
$ws.Range("U21").Value = "null"
$ws.Range("V21").Value = "lance"
$ws.Range("W21").Value = "null"
$ws.Range("X21").Value = "null"
$ws.Range("Y21").Value = "END"

# --- New row 22: "Illyia?" ---
$ws.Range("A22").Value = "Illyia?"
$ws.Range("B22").Value = "spr_bt_changeling"
$ws.Range("C22").Value = "enemy_general_1"
$ws.Range("D22").Value = 150
$ws.Range("E22").Value = 0
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 1
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0.5
$ws.Range("J22").Value = 1
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = "Attacker1"
$ws.Range("M22").Value = "none"
$ws.Range("N22").Value = -2
$ws.Range("O22").Value = -2
$ws.Range("P22").Value = -2
$ws.Range("Q22").Value = -2
$ws.Range("R22").Value = "GROUP/DEFAULT"
$ws.Range("S22").Value = "stress"
$ws.Range("T22").Value = "cryo"
$ws.Range("U22").Value = "other"
$ws.Range("V22").Value = "stress"
$ws.Range("W22").Value = "null"
$ws.Range("X22").Value = "other"
$ws.Range("Y22").Value = "stress"
$ws.Range("Z22").Value = "agni"
$ws.Range("AA22").Value = "other"
$ws.Range("AB22").Value = "END"

# --- New row 23: "Dahlia?" ---
$ws.Range("A23").Value = "Dahlia?"
$ws.Range("B23").Value = "spr_bt_changeling"
$ws.Range("C23").Value = "enemy_general_1"
$ws.Range("D23").Value = 150
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 2
$ws.Range("G23").Value = 1
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0.9
$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = "Attacker1"
$ws.Range("M23").Value = "none"
$ws.Range("N23").Value = -2
$ws.Range("O23").Value = -2
$ws.Range("P23").Value = -2
$ws.Range("Q23").Value = -2
$ws.Range("R23").Value = "GROUP/DEFAULT"
$ws.Range("S23").Value = "stress"
$ws.Range("T23").Value = "null"
$ws.Range("U23").Value = "other"
$ws.Range("V23").Value = "stress"
$ws.Range("W23").Value = "null"
$ws.Range("X23").Value = "other"
$ws.Range("Y23").Value = "stress"
$ws.Range("Z23").Value = "veld"
$ws.Range("AA23").Value = "other"
$ws.Range("AB23").Value = "END"

# Keep the cursor where the author's last save left it.
$ws.Range("I24").Select()
